$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ------------------------------------------------------------------
# 1) Insert a new column before N (shift old N:AC -> O:AD)
# ------------------------------------------------------------------
$ws.Range("N1").EntireColumn.Insert(-4161)

# The insert grows the <cols> customWidth span by one column; put it
# back the way it was (the new rightmost column carries no explicit
# per-cell styling anyway, so this is purely cosmetic bookkeeping).
$ws.Columns.Item(19).ClearFormats()

# ------------------------------------------------------------------
# 2) Populate the new column N with the "localdb" command group
# ------------------------------------------------------------------
$ws.Cells.Item(1,14).Value = "localdb"
$ws.Cells.Item(2,14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3,14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4,14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5,14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6,14).Value = "purge(var)"
$ws.Cells.Item(7,14).Value = "runSQLs(var,sqls)"

# ------------------------------------------------------------------
# 3) Shift column A (the "target" index list) down by one row
#    starting at row 14, so that "localdb" can be inserted
#    alphabetically between "json" and "macro". We use Copy
#    (instead of Range.Insert) because Range.Insert on this engine
#    shifts the *entire* row rather than just the target column.
# ------------------------------------------------------------------
$ws.Range("A14:A29").Copy($ws.Range("A15:A30"))
$ws.Cells.Item(14,1).Value = "localdb"

# ------------------------------------------------------------------
# 4) Fix up the defined names so they point at their new locations
# ------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo      = '=''#system''!$P$2:$P$2'
$wb.Names.Item("number").RefersTo    = '=''#system''!$Q$2:$Q$16'
$wb.Names.Item("pdf").RefersTo       = '=''#system''!$R$2:$R$16'
$wb.Names.Item("rdbms").RefersTo     = '=''#system''!$S$2:$S$7'
$wb.Names.Item("redis").RefersTo     = '=''#system''!$T$2:$T$10'
$wb.Names.Item("sms").RefersTo       = '=''#system''!$U$2:$U$2'
$wb.Names.Item("sound").RefersTo     = '=''#system''!$V$2:$V$5'
$wb.Names.Item("ssh").RefersTo       = '=''#system''!$W$2:$W$9'
$wb.Names.Item("step").RefersTo      = '=''#system''!$X$2:$X$4'
$wb.Names.Item("target").RefersTo    = '=''#system''!$A$2:$A$30'
$wb.Names.Item("web").RefersTo       = '=''#system''!$Y$2:$Y$127'
$wb.Names.Item("webalert").RefersTo  = '=''#system''!$Z$2:$Z$8'
$wb.Names.Item("webcookie").RefersTo = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("ws").RefersTo        = '=''#system''!$AB$2:$AB$17'
$wb.Names.Item("ws.async").RefersTo  = '=''#system''!$AC$2:$AC$8'
$wb.Names.Item("xml").RefersTo       = '=''#system''!$AD$2:$AD$21'
$wb.Names.Item("macro").RefersTo     = '=''#system''!$O$2:$O$4'

$wb.Names.Add("localdb", '=''#system''!$N$2:$N$7')
